$d = $word.ActiveDocument

# The heading paragraph reads:
#   "PROGRAM" + <13 spaces, sz=24> + <ECE 1ST SEM ...>
# We need to turn the 13-space run into:
#   " " (sz=26) + "-" (bold, sz=26) + "12" (bold, sz=26)
#   + " " (bold, sz=24) + "          " (10 spaces, sz=24)
# i.e. insert " -12 " (bold "-12") in the middle of the run of spaces,
# bumping the size of the leading part to 26 and bolding "-12 ".

$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "PROGRAM             "
$found = $find.Execute()
if (-not $found) {
    throw "Could not locate the 'PROGRAM' + 13-space anchor text"
}

# $rng now spans "PROGRAM" plus the 13 trailing spaces; the spaces are the
# last 13 characters of that match.
$spacesEnd = $rng.End
$spacesStart = $spacesEnd - 13

$target = $d.Range($spacesStart, $spacesEnd)
if ($target.Text -ne "             ") {
    throw "Unexpected text in target range: [$($target.Text)]"
}

# Replace the 13 spaces with the new 15-character sequence in one go, then
# reapply per-character formatting below (the new text first inherits the
# original run's formatting: not bold, sz=24/szCs=24).
$target.Text = " -12           "

# Run 1: " " -> sz 26
$run1 = $d.Range($spacesStart, $spacesStart + 1)
$run1.Font.Size = 13

# Run 2: "-" -> bold, sz 26
$run2 = $d.Range($spacesStart + 1, $spacesStart + 2)
$run2.Font.Bold = 1
$run2.Font.Size = 13

# Run 3: "12" -> bold, sz 26
$run3 = $d.Range($spacesStart + 2, $spacesStart + 4)
$run3.Font.Bold = 1
$run3.Font.Size = 13

# Run 4: " " -> bold, sz 24 (already sz 24 from the original run)
$run4 = $d.Range($spacesStart + 4, $spacesStart + 5)
$run4.Font.Bold = 1

# Run 5: 10 trailing spaces -> sz 24 (unchanged from the original run, no
# formatting change required).
